# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# 1) Rows 117 and 120 have their match data swapped (id/Div/Div name/
#    Date stay fixed per row, the rest of the columns exchange places).
# -------------------------------------------------------------------

# New values for row 117 (previously held by row 120)
$ws.Cells.Item(117, 2).Value = 7013409          # B - id/match number
$ws.Cells.Item(117, 6).Value = "Nacional De Football"   # F - HomeTeam
$ws.Cells.Item(117, 7).Value = "Torque"                 # G - AwayTeam
$ws.Cells.Item(117, 8).Value = 1                # H - FTHG
$ws.Cells.Item(117, 9).Value = 1                # I - FTAG
$ws.Cells.Item(117, 10).Value = "D"             # J - FTR
$ws.Cells.Item(117, 11).Value = 1.666           # K - oddH_op
$ws.Cells.Item(117, 12).Value = 3.9             # L - oddD_op
$ws.Cells.Item(117, 13).Value = 4.5             # M - oddA_op
$ws.Cells.Item(117, 14).Value = 1.615           # N - oddH
$ws.Cells.Item(117, 15).Value = 4               # O - oddD
$ws.Cells.Item(117, 16).Value = 4.75            # P - oddA
$ws.Cells.Item(117, 17).Value = -0.75           # Q - Ah
$ws.Cells.Item(117, 18).Value = 1.8             # R - oddAHH_op
$ws.Cells.Item(117, 19).Value = 2.05            # S - oddAHA_op
$ws.Cells.Item(117, 20).Value = 2.75            # T - AhOU_op
$ws.Cells.Item(117, 21).Value = 1.95            # U - oddAHH
$ws.Cells.Item(117, 22).Value = 1.9             # V - oddAHA
$ws.Cells.Item(117, 23).Value = -1              # W
$ws.Cells.Item(117, 24).Value = 3               # X
$ws.Cells.Item(117, 25).Value = -1              # Y
$ws.Cells.Item(117, 26).Value = -1              # Z
$ws.Cells.Item(117, 27).Value = 1.05            # AA
$ws.Cells.Item(117, 28).Value = -1              # AB
$ws.Cells.Item(117, 29).Value = 0.8999999999999999 # AC

# New values for row 120 (previously held by row 117)
$ws.Cells.Item(120, 2).Value = 7013886
$ws.Cells.Item(120, 6).Value = "Racing Club de Montevideo"
$ws.Cells.Item(120, 7).Value = "Cerro"
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 1
$ws.Cells.Item(120, 10).Value = "A"
$ws.Cells.Item(120, 11).Value = 2.25
$ws.Cells.Item(120, 12).Value = 3.1
$ws.Cells.Item(120, 13).Value = 3.25
$ws.Cells.Item(120, 14).Value = 2.25
$ws.Cells.Item(120, 15).Value = 2.875
$ws.Cells.Item(120, 16).Value = 3.5
$ws.Cells.Item(120, 17).Value = -0.25
$ws.Cells.Item(120, 18).Value = 1.95
$ws.Cells.Item(120, 19).Value = 1.9
$ws.Cells.Item(120, 20).Value = 2
$ws.Cells.Item(120, 21).Value = 1.925
$ws.Cells.Item(120, 22).Value = 1.925
$ws.Cells.Item(120, 23).Value = -1
$ws.Cells.Item(120, 24).Value = -1
$ws.Cells.Item(120, 25).Value = 2.5
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 0.8999999999999999
$ws.Cells.Item(120, 28).Value = -1
$ws.Cells.Item(120, 29).Value = 0.925

# -------------------------------------------------------------------
# 2) Append three brand new match rows (169, 170, 171) at the bottom
#    of the sheet, copying the formatting of row 168 for columns A
#    (bold/bordered id style) and E (date style).
# -------------------------------------------------------------------

$ws.Range("A168").Copy($ws.Range("A169:A171"))
$ws.Range("E168").Copy($ws.Range("E169:E171"))

# Row 169
$ws.Cells.Item(169, 1).Value = 167
$ws.Cells.Item(169, 2).Value = 8014131
$ws.Cells.Item(169, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(169, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(169, 5).Value = 45387.75
$ws.Cells.Item(169, 6).Value = "Boston River"
$ws.Cells.Item(169, 7).Value = "Defensor Sporting"
$ws.Cells.Item(169, 11).Value = 2.875
$ws.Cells.Item(169, 12).Value = 3.3
$ws.Cells.Item(169, 13).Value = 2.375
$ws.Cells.Item(169, 14).Value = 2.8
$ws.Cells.Item(169, 15).Value = 3.3
$ws.Cells.Item(169, 16).Value = 2.4
$ws.Cells.Item(169, 17).Value = 0
$ws.Cells.Item(169, 18).Value = 2.1
$ws.Cells.Item(169, 19).Value = 1.775
$ws.Cells.Item(169, 20).Value = 2.5
$ws.Cells.Item(169, 21).Value = 1.95
$ws.Cells.Item(169, 22).Value = 1.9
$ws.Cells.Item(169, 23).Value = 0
$ws.Cells.Item(169, 24).Value = 0
$ws.Cells.Item(169, 25).Value = 0
$ws.Cells.Item(169, 26).Value = 0
$ws.Cells.Item(169, 27).Value = 0

# Row 170
$ws.Cells.Item(170, 1).Value = 168
$ws.Cells.Item(170, 2).Value = 8014133
$ws.Cells.Item(170, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(170, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(170, 5).Value = 45388.41666666666
$ws.Cells.Item(170, 6).Value = "CA River Plate"
$ws.Cells.Item(170, 7).Value = "Montevideo Wanderers"
$ws.Cells.Item(170, 11).Value = 2.5
$ws.Cells.Item(170, 12).Value = 2.9
$ws.Cells.Item(170, 13).Value = 2.9
$ws.Cells.Item(170, 14).Value = 2.5
$ws.Cells.Item(170, 15).Value = 2.9
$ws.Cells.Item(170, 16).Value = 2.875
$ws.Cells.Item(170, 17).Value = 0
$ws.Cells.Item(170, 18).Value = 1.8
$ws.Cells.Item(170, 19).Value = 2.05
$ws.Cells.Item(170, 20).Value = 2
$ws.Cells.Item(170, 21).Value = 1.85
$ws.Cells.Item(170, 22).Value = 2
$ws.Cells.Item(170, 23).Value = 0
$ws.Cells.Item(170, 24).Value = 0
$ws.Cells.Item(170, 25).Value = 0
$ws.Cells.Item(170, 26).Value = 0
$ws.Cells.Item(170, 27).Value = 0

# Row 171
$ws.Cells.Item(171, 1).Value = 169
$ws.Cells.Item(171, 2).Value = 8014091
$ws.Cells.Item(171, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(171, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(171, 5).Value = 45389.70833333334
$ws.Cells.Item(171, 6).Value = "Deportivo Maldonado"
$ws.Cells.Item(171, 7).Value = "Penarol"
$ws.Cells.Item(171, 11).Value = 5
$ws.Cells.Item(171, 12).Value = 3.75
$ws.Cells.Item(171, 13).Value = 1.615
$ws.Cells.Item(171, 14).Value = 4.5
$ws.Cells.Item(171, 15).Value = 3.6
$ws.Cells.Item(171, 16).Value = 1.7
$ws.Cells.Item(171, 17).Value = 0.75
$ws.Cells.Item(171, 18).Value = 1.875
$ws.Cells.Item(171, 19).Value = 1.975
$ws.Cells.Item(171, 20).Value = 2.5
$ws.Cells.Item(171, 21).Value = 2.025
$ws.Cells.Item(171, 22).Value = 1.825
$ws.Cells.Item(171, 23).Value = 0
$ws.Cells.Item(171, 24).Value = 0
$ws.Cells.Item(171, 25).Value = 0
$ws.Cells.Item(171, 26).Value = 0
$ws.Cells.Item(171, 27).Value = 0
